# Presupuesto De Costes.xlsx - "Cambios en la Docu"
#
# The project's hourly rate drops (N8: 7.5 -> 5), one material cost entry is
# renegotiated (F15: 500 -> 150), a chunk of the hour estimates for the
# "Desarrollo" phase are reduced, the profit margin drops from 14% to 7%,
# the sub-phase previously labelled "Fase de Desarrollo" at D43 is relabelled
# "Fase de Pruebas", and the client-facing budget figures (S61, S62, S64:S70)
# are updated to match. All dependent subtotals/totals are plain formulas
# that recompute automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Hourly rate used throughout the "Horas" columns (I = H * $N$8) -------
$ws.Range("N8").Value = 5

# --- Material cost line (Capítulo I.13, "Fase de Desarrollo" materials) ---
$ws.Range("F15").Value = 150

# --- Hours per task, Capítulo C.2 "Fase de Desarrollo" (rows 22-50) -------
$ws.Range("H22").Value = 4
$ws.Range("H24").Value = 3
$ws.Range("H25").Value = 3
$ws.Range("H28").Value = 3
$ws.Range("H29").Value = 8
$ws.Range("H30").Value = 8
$ws.Range("H31").Value = 8
$ws.Range("H33").Value = 8
$ws.Range("H34").Value = 6
$ws.Range("H35").Value = 4
$ws.Range("H38").Value = 18
$ws.Range("H39").Value = 18
$ws.Range("H40").Value = 18
$ws.Range("H41").Value = 12
$ws.Range("H44").Value = 8
$ws.Range("H45").Value = 18
$ws.Range("H48").Value = 18
$ws.Range("H49").Value = 8
$ws.Range("H50").Value = 3

# --- Relabel the sub-phase header at D43 (was "Fase de Desarrollo") -------
# Set the BENEFICIO label first so it reuses its existing shared-string slot
# in place, then add the new "Fase de Pruebas" label.
$ws.Range("A52").Value = "BENEFICIO (7%):"
$ws.Range("D43").Value = "Fase de Pruebas"

# --- Profit margin formula: 14% -> 7% --------------------------------------
$ws.Range("J52").Formula = "=J51*0.07"

# --- Client budget figures (right-hand "PRESUPUESTO DEL CLIENTE" table) ---
$ws.Range("S61").Value = 285
$ws.Range("S62").Value = 90
$ws.Range("S64").Value = 132.5
$ws.Range("S65").Value = 132
$ws.Range("S66").Value = 132.2
$ws.Range("S67").Value = 92.04
$ws.Range("S68").Value = 72
$ws.Range("S69").Value = 190
$ws.Range("S70").Value = 145

# --- View state: re-centre the window, zoom out, move the selection -------
$excel.ActiveWindow.ScrollRow = 26
$excel.ActiveWindow.ScrollColumn = 1
$excel.ActiveWindow.Zoom = 70
$ws.Range("S68").Select()
